$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the 380 kV case (pl_mw.xlsx row data, columns B:O except F,I,K,L)
$dataRowCount = 24

$colB = New-Object 'object[,]' $dataRowCount,1
$colB[0,0] = 0.7598684979915902
$colB[1,0] = 0.6830668096117449
$colB[2,0] = 0.6359712845096794
$colB[3,0] = 0.6167959028718428
$colB[4,0] = 0.6136128694565173
$colB[5,0] = 0.6357126106570661
$colB[6,0] = 0.7333752357078538
$colB[7,0] = 0.9253375082166144
$colB[8,0] = 1.06660715973544
$colB[9,0] = 1.130918428076882
$colB[10,0] = 1.155277284821409
$colB[11,0] = 1.15003094003572
$colB[12,0] = 1.132922340913694
$colB[13,0] = 1.122443520938475
$colB[14,0] = 1.06240511161559
$colB[15,0] = 1.025584686527054
$colB[16,0] = 1.004411045850361
$colB[17,0] = 0.9972428228508079
$colB[18,0] = 1.029503827100712
$colB[19,0] = 1.137947409300637
$colB[20,0] = 1.208853508216805
$colB[21,0] = 1.171007055629843
$colB[22,0] = 1.027731999405148
$colB[23,0] = 0.8733626994046517
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' $dataRowCount,1
$colC[0,0] = 0.126138888061746
$colC[1,0] = 0.1100499609132441
$colC[2,0] = 0.1001338299085717
$colC[3,0] = 0.0960837463780706
$colC[4,0] = 0.09541068577860301
$colC[5,0] = 0.1000792459127524
$colC[6,0] = 0.1205993415152875
$colC[7,0] = 0.1605332070531915
$colC[8,0] = 0.189677145595283
$colC[9,0] = 0.2028913110076189
$colC[10,0] = 0.207888700609459
$colC[11,0] = 0.2068127178841053
$colC[12,0] = 0.2033025816008944
$colC[13,0] = 0.2011516637954571
$colC[14,0] = 0.1888126696821359
$colC[15,0] = 0.1812317517316728
$colC[16,0] = 0.1768673240717646
$colC[17,0] = 0.1753889104659265
$colC[18,0] = 0.1820391780028672
$colC[19,0] = 0.2043337729876953
$colC[20,0] = 0.2188664069418564
$colC[21,0] = 0.2111136466679682
$colC[22,0] = 0.1816741594348628
$colC[23,0] = 0.1497637225193671
$ws.Range("C2:C25").Value = $colC

$colD = New-Object 'object[,]' $dataRowCount,1
$colD[0,0] = 0.6026544702145884
$colD[1,0] = 0.5985411035631074
$colD[2,0] = 0.5963047227232607
$colD[3,0] = 0.5954662074588839
$colD[4,0] = 0.5953313741008657
$colD[5,0] = 0.5962931192067771
$colD[6,0] = 0.6011762035574719
$colD[7,0] = 0.6130439380945631
$colD[8,0] = 0.6231582072106789
$colD[9,0] = 0.6280620178646643
$colD[10,0] = 0.6299624366196781
$colD[11,0] = 0.629551216282124
$colD[12,0] = 0.6282174961147575
$colD[13,0] = 0.6274062101594211
$colD[14,0] = 0.6228438168858474
$colD[15,0] = 0.6201224232980849
$colD[16,0] = 0.6185856538757548
$colD[17,0] = 0.6180702284347888
$colD[18,0] = 0.6204091708467843
$colD[19,0] = 0.6286080634743314
$colD[20,0] = 0.6342197453884921
$colD[21,0] = 0.6312015432050089
$colD[22,0] = 0.6202794456924892
$colD[23,0] = 0.6095882562475907
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'object[,]' $dataRowCount,1
$colE[0,0] = 0.2310343479695582
$colE[1,0] = 0.2306167796436718
$colE[2,0] = 0.2304693245033569
$colE[3,0] = 0.2304366493051475
$colE[4,0] = 0.230432879994666
$colE[5,0] = 0.2304687728069013
$colE[6,0] = 0.2308677735849898
$colE[7,0] = 0.2325139156336462
$colE[8,0] = 0.2342494923505356
$colE[9,0] = 0.2351532546132589
$colE[10,0] = 0.2355119005966202
$colE[11,0] = 0.2354339301562902
$colE[12,0] = 0.2351824318320794
$colE[13,0] = 0.2350305185265391
$colE[14,0] = 0.2341927267268709
$colE[15,0] = 0.2337080162171681
$colE[16,0] = 0.2334399773290379
$colE[17,0] = 0.2333510714486877
$colE[18,0] = 0.2337585016242478
$colE[19,0] = 0.2352558577339963
$colE[20,0] = 0.2363301136789246
$colE[21,0] = 0.235748016840958
$colE[22,0] = 0.233735644066666
$colE[23,0] = 0.2319761649496108
$ws.Range("E2:E25").Value = $colE

$colG = New-Object 'object[,]' $dataRowCount,1
$colG[0,0] = 0.6688545619091002
$colG[1,0] = 0.6663194885166064
$colG[2,0] = 0.6652471898419634
$colG[3,0] = 0.6649318155363204
$colG[4,0] = 0.6648867869616026
$colG[5,0] = 0.6652424444846758
$colG[6,0] = 0.6678798570898152
$colG[7,0] = 0.6769031894488506
$colG[8,0] = 0.6858950695130233
$colG[9,0] = 0.6905019327406166
$colG[10,0] = 0.6923209143996871
$colG[11,0] = 0.6919258491859637
$colG[12,0] = 0.6906500879556887
$colG[13,0] = 0.6898783506003241
$colG[14,0] = 0.6856044090774418
$colG[15,0] = 0.683114894656029
$colG[16,0] = 0.6817315868409821
$colG[17,0] = 0.6812715626264634
$colG[18,0] = 0.6833748768374051
$colG[19,0] = 0.6910227874472383
$colG[20,0] = 0.69645525748426
$colG[21,0] = 0.6935160572231212
$colG[22,0] = 0.6832571895554764
$colG[23,0] = 0.674048318034778
$ws.Range("G2:G25").Value = $colG

$colH = New-Object 'object[,]' $dataRowCount,1
$colH[0,0] = 0.7666454822131499
$colH[1,0] = 0.769877040480651
$colH[2,0] = 0.7722405516814632
$colH[3,0] = 0.7732991052117484
$colH[4,0] = 0.7734806399636796
$colH[5,0] = 0.7722544413944803
$colH[6,0] = 0.7676810062564243
$colH[7,0] = 0.7617219289796253
$colH[8,0] = 0.7591787938745398
$colH[9,0] = 0.7584204609402292
$colH[10,0] = 0.7581906128366285
$colH[11,0] = 0.758237565507585
$colH[12,0] = 0.7584004024689506
$colH[13,0] = 0.7585076092103833
$colH[14,0] = 0.7592363730052085
$colH[15,0] = 0.7597855280345129
$colH[16,0] = 0.7601389007636499
$colH[17,0] = 0.7602649895010671
$colH[18,0] = 0.759723187025088
$colH[19,0] = 0.7583510177225321
$colH[20,0] = 0.7577883104838037
$colH[21,0] = 0.7580580674632103
$colH[22,0] = 0.759751254096841
$colH[23,0] = 0.7630118013194789
$ws.Range("H2:H25").Value = $colH

$colJ = New-Object 'object[,]' $dataRowCount,1
$colJ[0,0] = 0.1103907924338472
$colJ[1,0] = 0.1108896302811608
$colJ[2,0] = 0.1112448835794133
$colJ[3,0] = 0.1114019670293054
$colJ[4,0] = 0.1114287944711663
$colJ[5,0] = 0.1112469521965949
$colJ[6,0] = 0.1105526302825623
$colJ[7,0] = 0.1095795453920729
$colJ[8,0] = 0.1091014917759239
$colJ[9,0] = 0.1089354658488269
$colJ[10,0] = 0.1088799941015068
$colJ[11,0] = 0.1088916118565741
$colJ[12,0] = 0.1089307538671065
$colJ[13,0] = 0.1089556930413558
$colJ[14,0] = 0.1091133772291784
$colJ[15,0] = 0.1092232880136734
$colJ[16,0] = 0.1092913477771518
$colJ[17,0] = 0.1093152232237848
$colJ[18,0] = 0.1092110866872851
$colJ[19,0] = 0.1089190561023656
$colJ[20,0] = 0.1087713229542331
$colJ[21,0] = 0.108846224485788
$colJ[22,0] = 0.109216587735915
$colJ[23,0] = 0.1098011970450123
$ws.Range("J2:J25").Value = $colJ

$colM = New-Object 'object[,]' $dataRowCount,1
$colM[0,0] = 0.3881540985016514
$colM[1,0] = 0.3660828348957423
$colM[2,0] = 0.3526610250504802
$colM[3,0] = 0.3472245294209699
$colM[4,0] = 0.3463238055566649
$colM[5,0] = 0.3525875725484724
$colM[6,0] = 0.3805171118537274
$colM[7,0] = 0.4363083525099114
$colM[8,0] = 0.4779115540984833
$colM[9,0] = 0.4969694107164244
$colM[10,0] = 0.5042049138018569
$colM[11,0] = 0.5026457919605107
$colM[12,0] = 0.4975643068693358
$colM[13,0] = 0.4944541769946866
$colM[14,0] = 0.4766687137572561
$colM[15,0] = 0.4657915747336361
$colM[16,0] = 0.4595478090476277
$colM[17,0] = 0.4574359301323909
$colM[18,0] = 0.4669481757899163
$colM[19,0] = 0.4990563571242745
$colM[20,0] = 0.5201497941201012
$colM[21,0] = 0.5088819743581681
$colM[22,0] = 0.4664252463808012
$colM[23,0] = 0.4211069543555226
$ws.Range("M2:M25").Value = $colM

$colN = New-Object 'object[,]' $dataRowCount,1
$colN[0,0] = 1.245465653723244
$colN[1,0] = 1.25709921233765
$colN[2,0] = 1.264733556251187
$colN[3,0] = 1.267968216742688
$colN[4,0] = 1.26851279695321
$colN[5,0] = 1.264776679462649
$colN[6,0] = 1.249374975063674
$colN[7,0] = 1.223067626624392
$colN[8,0] = 1.206110483037698
$colN[9,0] = 1.198910325060929
$colN[10,0] = 1.196257649939454
$colN[11,0] = 1.196825665711593
$colN[12,0] = 1.198690607464769
$colN[13,0] = 1.199842557899778
$colN[14,0] = 1.206591367313543
$colN[15,0] = 1.210863116947777
$colN[16,0] = 1.213368469252771
$colN[17,0] = 1.214225043741024
$colN[18,0] = 1.21040337740834
$colN[19,0] = 1.198140824118973
$colN[20,0] = 1.190557099018257
$colN[21,0] = 1.194565280040223
$colN[22,0] = 1.210611071565459
$colN[23,0] = 1.229767852702665
$ws.Range("N2:N25").Value = $colN

$colO = New-Object 'object[,]' $dataRowCount,1
$colO[0,0] = 2.864803467933086
$colO[1,0] = 2.865872919563913
$colO[2,0] = 2.868338661141337
$colO[3,0] = 2.869798052373909
$colO[4,0] = 2.870067829177742
$colO[5,0] = 2.868356502861445
$colO[6,0] = 2.864796530407517
$colO[7,0] = 2.872188738342459
$colO[8,0] = 2.886414262141869
$colO[9,0] = 2.894802758947208
$colO[10,0] = 2.898255444685077
$colO[11,0] = 2.897499558329343
$colO[12,0] = 2.895081276580129
$colO[13,0] = 2.893635984568448
$colO[14,0] = 2.885904672133989
$colO[15,0] = 2.881653131290506
$colO[16,0] = 2.879388189907814
$colO[17,0] = 2.878652295462729
$colO[18,0] = 2.882087037734181
$colO[19,0] = 2.895784086513657
$colO[20,0] = 2.906345685766155
$colO[21,0] = 2.900561318175761
$colO[22,0] = 2.88189031006857
$colO[23,0] = 2.868646868223635
$ws.Range("O2:O25").Value = $colO
